$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename the "Include from ActMood" sheet to "Include #0" ---
$ws2.Name = "Include #0"

# --- Update Version metadata value (row 3, column B) ---
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# --- Update Date metadata value (row 8, column B) ---
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Insert a new "Jurisdiction" row after "Contact" (row 10), pushing the
#     Description / Purpose / Copyright / Immutable rows down by one ---

# Read the existing values (top to bottom) before we start overwriting cells,
# so the shift doesn't clobber data we still need to move.
$v14a = $ws1.Range("A14").Value2
$v14b = $ws1.Range("B14").Value2
$v13a = $ws1.Range("A13").Value2
$v13b = $ws1.Range("B13").Value2
$v12a = $ws1.Range("A12").Value2
$v12b = $ws1.Range("B12").Value2
$v11a = $ws1.Range("A11").Value2
$v11b = $ws1.Range("B11").Value2

# Give the new last row (15) the same formatting as the other body rows.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)

# Shift content down one row.
$ws1.Range("A15").Value = $v14a
$ws1.Range("B15").Value = $v14b
$ws1.Range("A14").Value = $v13a
$ws1.Range("B14").Value = $v13b
$ws1.Range("A13").Value = $v12a
$ws1.Range("B13").Value = $v12b
$ws1.Range("A12").Value = $v11a
$ws1.Range("B12").Value = $v11b

# Row 11 becomes the new "Jurisdiction" property with a blank (empty-string) value.
$ws1.Range("A11").Value = "Jurisdiction"
$ws2.Range("B9").Copy()
$ws1.Range("B11").PasteSpecial(-4163)

Write-Host "Edit applied"
